$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02072366666666666
$ws.Range("H2").Value = 0.062171
$ws.Range("I2").Value = 0.03659418775688948
$ws.Range("J2").Value = 0.03659418775688947
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1311436666666667
$ws.Range("N2").Value = 0.393431
$ws.Range("O2").Value = 0.02663441993971509
$ws.Range("P2").Value = 0.02663441993971509
$ws.Range("Q2").Value = 0.002717777633444444
$ws.Range("R2").Value = 0.024459998701
$ws.Range("S2").Value = 0.0009746649640697751
$ws.Range("T2").Value = 0.0009746649640697748

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02072366666666666
$ws.Range("H3").Value = 0.062171
$ws.Range("I3").Value = 0.03659418775688948
$ws.Range("J3").Value = 0.03659418775688947
$ws.Range("O3").Value = 0.06149297381279183
$ws.Range("P3").Value = 0.06149297381279183
$ws.Range("Q3").Value = 0.006274746332777777
$ws.Range("R3").Value = 0.05647271699499999
$ws.Range("S3").Value = 0.002250285429434792
$ws.Range("T3").Value = 0.002250285429434792

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02072366666666666
$ws.Range("H4").Value = 0.062171
$ws.Range("I4").Value = 0.03659418775688948
$ws.Range("J4").Value = 0.03659418775688947
$ws.Range("M4").Value = 4.009307333333333
$ws.Range("N4").Value = 12.027922
$ws.Range("O4").Value = 0.8142640654908683
$ws.Range("P4").Value = 0.8142640654908684
$ws.Range("Q4").Value = 0.08308754874022221
$ws.Range("R4").Value = 0.747787938662
$ws.Range("S4").Value = 0.02979733209626099
$ws.Range("T4").Value = 0.02979733209626098

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02072366666666666
$ws.Range("H5").Value = 0.062171
$ws.Range("I5").Value = 0.03659418775688948
$ws.Range("J5").Value = 0.03659418775688947
$ws.Range("M5").Value = 0.480609
$ws.Range("N5").Value = 1.441827
$ws.Range("O5").Value = 0.09760854075662465
$ws.Range("P5").Value = 0.09760854075662465
$ws.Range("Q5").Value = 0.009959980713
$ws.Range("R5").Value = 0.089639826417
$ws.Range("S5").Value = 0.003571905267123922
$ws.Range("T5").Value = 0.003571905267123921

$ws.Range("I6").Value = 0.4063843675817323
$ws.Range("J6").Value = 0.4063843675817323
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1311436666666667
$ws.Range("N6").Value = 0.393431
$ws.Range("O6").Value = 0.02663441993971509
$ws.Range("P6").Value = 0.02663441993971509
$ws.Range("Q6").Value = 0.03018135973211111
$ws.Range("R6").Value = 0.271632237589
$ws.Range("S6").Value = 0.0108238119031074
$ws.Range("T6").Value = 0.0108238119031074

$ws.Range("I7").Value = 0.4063843675817323
$ws.Range("J7").Value = 0.4063843675817323
$ws.Range("O7").Value = 0.06149297381279183
$ws.Range("P7").Value = 0.06149297381279183
$ws.Range("S7").Value = 0.02498978327363144
$ws.Range("T7").Value = 0.02498978327363144

$ws.Range("I8").Value = 0.4063843675817323
$ws.Range("J8").Value = 0.4063843675817323
$ws.Range("M8").Value = 4.009307333333333
$ws.Range("N8").Value = 12.027922
$ws.Range("O8").Value = 0.8142640654908683
$ws.Range("P8").Value = 0.8142640654908684
$ws.Range("Q8").Value = 0.9227006532575555
$ws.Range("R8").Value = 8.304305879318001
$ws.Range("S8").Value = 0.3309041872990368
$ws.Range("T8").Value = 0.3309041872990368

$ws.Range("I9").Value = 0.4063843675817323
$ws.Range("J9").Value = 0.4063843675817323
$ws.Range("M9").Value = 0.480609
$ws.Range("N9").Value = 1.441827
$ws.Range("O9").Value = 0.09760854075662465
$ws.Range("P9").Value = 0.09760854075662465
$ws.Range("Q9").Value = 0.110607195057
$ws.Range("R9").Value = 0.995464755513
$ws.Range("S9").Value = 0.03966658510595665
$ws.Range("T9").Value = 0.03966658510595665

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.05227533333333333
$ws.Range("H10").Value = 0.156826
$ws.Range("I10").Value = 0.09230863407636922
$ws.Range("J10").Value = 0.0923086340763692
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1311436666666667
$ws.Range("N10").Value = 0.393431
$ws.Range("O10").Value = 0.02663441993971509
$ws.Range("P10").Value = 0.02663441993971509
$ws.Range("Q10").Value = 0.006855578889555555
$ws.Range("R10").Value = 0.06170021000599999
$ws.Range("S10").Value = 0.002458586924051512
$ws.Range("T10").Value = 0.002458586924051512

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.05227533333333333
$ws.Range("H11").Value = 0.156826
$ws.Range("I11").Value = 0.09230863407636922
$ws.Range("J11").Value = 0.0923086340763692
$ws.Range("O11").Value = 0.06149297381279183
$ws.Range("P11").Value = 0.06149297381279183
$ws.Range("Q11").Value = 0.01582801255222222
$ws.Range("R11").Value = 0.14245211297
$ws.Range("S11").Value = 0.005676332417952756
$ws.Range("T11").Value = 0.005676332417952755

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.05227533333333333
$ws.Range("H12").Value = 0.156826
$ws.Range("I12").Value = 0.09230863407636922
$ws.Range("J12").Value = 0.0923086340763692
$ws.Range("M12").Value = 4.009307333333333
$ws.Range("N12").Value = 12.027922
$ws.Range("O12").Value = 0.8142640654908683
$ws.Range("P12").Value = 0.8142640654908684
$ws.Range("Q12").Value = 0.2095878772857778
$ws.Range("R12").Value = 1.886290895572
$ws.Range("S12").Value = 0.0751636036629333
$ws.Range("T12").Value = 0.0751636036629333

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.05227533333333333
$ws.Range("H13").Value = 0.156826
$ws.Range("I13").Value = 0.09230863407636922
$ws.Range("J13").Value = 0.0923086340763692
$ws.Range("M13").Value = 0.480609
$ws.Range("N13").Value = 1.441827
$ws.Range("O13").Value = 0.09760854075662465
$ws.Range("P13").Value = 0.09760854075662465
$ws.Range("Q13").Value = 0.025123995678
$ws.Range("R13").Value = 0.226115961102
$ws.Range("S13").Value = 0.009010111071431636
$ws.Range("T13").Value = 0.009010111071431634

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2631716666666666
$ws.Range("H14").Value = 0.789515
$ws.Range("I14").Value = 0.4647128105850091
$ws.Range("J14").Value = 0.464712810585009
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1311436666666667
$ws.Range("N14").Value = 0.393431
$ws.Range("O14").Value = 0.02663441993971509
$ws.Range("P14").Value = 0.02663441993971509
$ws.Range("Q14").Value = 0.03451329732944444
$ws.Range("R14").Value = 0.310619675965
$ws.Range("S14").Value = 0.01237735614848641
$ws.Range("T14").Value = 0.01237735614848641

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2631716666666666
$ws.Range("H15").Value = 0.789515
$ws.Range("I15").Value = 0.4647128105850091
$ws.Range("J15").Value = 0.464712810585009
$ws.Range("O15").Value = 0.06149297381279183
$ws.Range("P15").Value = 0.06149297381279183
$ws.Range("Q15").Value = 0.07968355585277777
$ws.Range("R15").Value = 0.7171520026749999
$ws.Range("S15").Value = 0.02857657269177286
$ws.Range("T15").Value = 0.02857657269177285

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2631716666666666
$ws.Range("H16").Value = 0.789515
$ws.Range("I16").Value = 0.4647128105850091
$ws.Range("J16").Value = 0.464712810585009
$ws.Range("M16").Value = 4.009307333333333
$ws.Range("N16").Value = 12.027922
$ws.Range("O16").Value = 0.8142640654908683
$ws.Range("P16").Value = 0.8142640654908684
$ws.Range("Q16").Value = 1.055136093092222
$ws.Range("R16").Value = 9.496224837829999
$ws.Range("S16").Value = 0.3783989424326373
$ws.Range("T16").Value = 0.3783989424326373

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2631716666666666
$ws.Range("H17").Value = 0.789515
$ws.Range("I17").Value = 0.4647128105850091
$ws.Range("J17").Value = 0.464712810585009
$ws.Range("M17").Value = 0.480609
$ws.Range("N17").Value = 1.441827
$ws.Range("O17").Value = 0.09760854075662465
$ws.Range("P17").Value = 0.09760854075662465
$ws.Range("Q17").Value = 0.126482671545
$ws.Range("R17").Value = 1.138344043905
$ws.Range("S17").Value = 0.04535993931211245
$ws.Range("T17").Value = 0.04535993931211244
